# The "added AAS lab details" commit re-saved the deck a few days later,
# which made PowerPoint refresh the cached text of the auto-updating
# "datetimeFigureOut" date field (footer date) shown on the slide master
# and on every slide layout, from 13-01-2023 to 16-01-2023.
#
# Slide.Master.CustomLayouts aliases every layout to the same object in
# this host, so CustomLayouts must be reached through
# Presentation.Designs(1).SlideMaster instead, which correctly yields
# each distinct layout.

$p = $ppt.ActivePresentation
$newDate = "16-01-2023"
$ppPlaceholderDate = 16

$slideMaster = $p.Designs.Item(1).SlideMaster
$layouts = $slideMaster.CustomLayouts

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($j = 1; $j -le $slideMaster.Shapes.Count; $j++) {
    $shape = $slideMaster.Shapes.Item($j)
    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}
